# Natmi follow-up: refresh LR-pair stats for Sema3a-Plxna4 to include the new "ECs" population
# (per Dr Hou's advice), expanding the sending/target cluster grid from 2x2 (FAPs, sCs) to 3x3
# (ECs, FAPs, sCs) and refreshing every computed column for all 9 row combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous 2x2 (4-row) result block before laying down the new 3x3 (9-row) block.
$ws.Range("A1:T10").ClearContents()

$headers = @(
    "Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Each inner array is one data row, columns A..T in order.
$rows = @(
    @("ECs", "Sema3a", "Plxna4", "ECs", 2, 0.6666666666666666, 0.5674196666666667, 1.702259, 0.07864125446886469, 0.07864125446886468, 2, 0.6666666666666666, 2.431926, 7.295778, 0.8146397085966534, 0.8146397085966534, 1.379922640278, 12.419303762502, 0.0640642886241912, 0.06406428862419118),
    @("ECs", "Sema3a", "Plxna4", "FAPs", 2, 0.6666666666666666, 0.5674196666666667, 1.702259, 0.07864125446886469, 0.07864125446886468, 3, 1, 0.220323, 0.6609689999999999, 0.07380317678956531, 0.07380317678956533, 0.125015603219, 1.125140428971, 0.005803974406518814, 0.005803974406518814),
    @("ECs", "Sema3a", "Plxna4", "sCs", 2, 0.6666666666666666, 0.5674196666666667, 1.702259, 0.07864125446886469, 0.07864125446886468, 3, 1, 0.333029, 0.999087, 0.1115571146137813, 0.1115571146137814, 0.1889672041703333, 1.700704837533, 0.008772991438154683, 0.008772991438154682),
    @("FAPs", "Sema3a", "Plxna4", "ECs", 3, 1, 0.7227763333333334, 2.168329, 0.1001728366019618, 0.1001728366019618, 2, 0.6666666666666666, 2.431926, 7.295778, 0.8146397085966534, 0.8146397085966534, 1.757738557218, 15.819647014962, 0.08160477041872233, 0.08160477041872233),
    @("FAPs", "Sema3a", "Plxna4", "FAPs", 3, 1, 0.7227763333333334, 2.168329, 0.1001728366019618, 0.1001728366019618, 3, 1, 0.220323, 0.6609689999999999, 0.07380317678956531, 0.07380317678956533, 0.159244250089, 1.433198250801, 0.007393073569246826, 0.007393073569246828),
    @("FAPs", "Sema3a", "Plxna4", "sCs", 3, 1, 0.7227763333333334, 2.168329, 0.1001728366019618, 0.1001728366019618, 3, 1, 0.333029, 0.999087, 0.1115571146137813, 0.1115571146137814, 0.2407054795136666, 2.166349315623, 0.01117499261399264, 0.01117499261399264),
    @("sCs", "Sema3a", "Plxna4", "ECs", 3, 1, 5.925096666666668, 17.77529, 0.8211859089291735, 0.8211859089291734, 2, 0.6666666666666666, 2.431926, 7.295778, 0.8146397085966534, 0.8146397085966534, 14.40939663618, 129.68456972562, 0.6689706495537399, 0.6689706495537397),
    @("sCs", "Sema3a", "Plxna4", "FAPs", 3, 1, 5.925096666666668, 17.77529, 0.8211859089291735, 0.8211859089291734, 3, 1, 0.220323, 0.6609689999999999, 0.07380317678956531, 0.07380317678956533, 1.30543507289, 11.74891565601, 0.06060612881379968, 0.06060612881379968),
    @("sCs", "Sema3a", "Plxna4", "sCs", 3, 1, 5.925096666666668, 17.77529, 0.8211859089291735, 0.8211859089291734, 3, 1, 0.333029, 0.999087, 0.1115571146137813, 0.1115571146137814, 1.973229017803333, 17.75906116023, 0.09160913056163403, 0.09160913056163401)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
